# Update the betting odds in row 2 (FlashScore weekly games export)
# on the active worksheet to reflect the latest scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.4
    "I2"  = 8
    "K2"  = 2.3
    "L2"  = 7.5
    "M2"  = 1.05
    "N2"  = 11
    "O2"  = 1.3
    "P2"  = 3.5
    "Q2"  = 1.98
    "R2"  = 1.92
    "U2"  = 2.2
    "V2"  = 1.62
    "W2"  = 6
    "X2"  = 6
    "AC2" = 10
    "AE2" = 23
    "AF2" = 81
    "AI2" = 21
    "AK2" = 51
    "AL2" = 51
    "AN2" = 3.2
    "AS2" = 151
    "AV2" = 67
    "AY2" = 41
    "AZ2" = 201
    "BA2" = 201
    "BB2" = 451
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
